# Insert one new data row at row 702 (pushing the existing rows 702-743
# down to 703-744) and populate it with the new entry:
#   2026/01/24, 土, 7, 201
#
# The sheet's date column (A) stores dates as literal text (e.g. "2026/01/24"),
# not as real Excel date serials, so the new value is written with a leading
# apostrophe to force text, then the auto-applied "quote prefix" style is
# cleared by resetting the cell style back to Normal so the cell matches the
# unstyled look of every other data cell in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(702).Insert()

$ws.Range("A702").Value = "'2026/01/24"
$ws.Range("A702").Style = "Normal"
$ws.Range("B702").Value = "土"
$ws.Range("C702").Value = 7
$ws.Range("D702").Value = 201
